# Adding test case to search OPQA-1244
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# ---------------------------------------------------------------------
# 1) Re-copy the formatting (styles) of row 113 onto new row 116 BEFORE
#    we touch row 113's own formatting, so row 116 ends up with the
#    A=7,B=7,C=3,D=7,E=7 style pattern that row 113 currently has.
# ---------------------------------------------------------------------
$ws.Range("A113:E113").Copy()
$ws.Range("A116:E116").PasteSpecial(-4122)  # xlPasteFormats

# ---------------------------------------------------------------------
# 2) Fill in the values for the new row 116 (new test case TestCase_B115)
#    Order matters for shared-string table append order: the new unique
#    strings must land at indices 362 (description), 363 (OPQA-1244),
#    364 (TestCase_B115) - so write C, then B, then A.
# ---------------------------------------------------------------------
$ws.Range("D116").Value = "Y"
$ws.Range("E116").Value = "SKIP"
$ws.Range("C116").Value = "Verify that search drop down content type is retained when user navigates back to PEOPLE search results page from profile page"
$ws.Range("B116").Value = "OPQA-1244"
$ws.Range("A116").Value = "TestCase_B115"

# ---------------------------------------------------------------------
# 3) Style-only fix ups on rows 112 and 113 (column D): s="7" -> s="3"
#    Copy the format from D110 (which already has style s="3").
# ---------------------------------------------------------------------
$ws.Range("D110").Copy()
$ws.Range("D112").PasteSpecial(-4122)
$ws.Range("D110").Copy()
$ws.Range("D113").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 4) Style fix ups on rows 114 and 115 (A:E): unstyled -> s="3"
#    Copy the format from row 110 (A:E all s="3").
# ---------------------------------------------------------------------
$ws.Range("A110:E110").Copy()
$ws.Range("A114:E114").PasteSpecial(-4122)
$ws.Range("A110:E110").Copy()
$ws.Range("A115:E115").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 5) Value changes
# ---------------------------------------------------------------------
$ws.Range("E107").Value = "SKIP"
$ws.Range("E109").Value = "SKIP"
$ws.Range("E110").Value = "SKIP"
$ws.Range("E111").Value = "SKIP"
$ws.Range("E113").Value = "PASS"
$ws.Range("E114").Value = "SKIP"
$ws.Range("E115").Value = "SKIP"

# ---------------------------------------------------------------------
# 6) Update the view: selected cell moves to the new last row (A116).
# ---------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 108
$win.ScrollColumn = 1
$ws.Range("A116").Select()
